$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.073.11"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.06%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.837.76"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.51%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9996"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.36"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.58%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6337"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.36%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07594"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +3.38%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2953"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.31%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.88"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.81%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07752"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.85%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.835.42"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.54%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.009"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.90%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6723"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "83.38"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.85%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009855"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +9.51%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.129"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.67%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.094.07"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.19%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.57"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.73%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "227.45"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.80%  "
$ws.Range("E21").Value = "  +0.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.258"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.77%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("D23").ClearFormats()
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "160.53"
$ws.Range("D24").ClearFormats()
$ws.Range("E25").Value = "  +4.16%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.561"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.54%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.00"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.08%  "
$ws.Range("E28").Value = "  +0.70%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.132"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.96%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.053"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.28%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.205"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.54%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.05388"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +2.32%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.867"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.52%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7535"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +3.08%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.146"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.32%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.671"
$ws.Range("D36").ClearFormats()
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.253.10"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -3.64%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01799"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.91%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.762"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.49%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.594"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +4.73%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9081"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.64%  "
$ws.Range("E42").Value = "  +0.29%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "102.84"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.94%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.980.74"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.45%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000124"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.27%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "65.12"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.92%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5118"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.08%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4101"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +3.41%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.054"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +2.57%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05808"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.34%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.791"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.93%  "
